$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.398.69'
$ws.Range("E2").Value = '  -7.58%  '
$ws.Range("D3").Value = '1.684.71'
$ws.Range("E3").Value = '  -6.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("E5").Value = '  -5.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5070'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -14.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.005'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2674'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '22.03'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06290'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07369'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.38%  '
$ws.Range("D12").Value = '1.687.38'
$ws.Range("E12").Value = '  -5.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.539'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5784'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.64%  '
$ws.Range("D15").Value = '1.913.73'
$ws.Range("E15").Value = '  -6.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008576'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.18'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -13.89%  '
$ws.Range("D18").Value = '26.453.18'
$ws.Range("E18").Value = '  -7.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.999'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.88%  '
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '185.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -11.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.257'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.006'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.63'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.479'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1166'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.349'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05725'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.336'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.522'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.515'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.656'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.17%  '
$ws.Range("E35").Value = '  -3.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.5947'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.356'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.661'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("D39").Value = '1.100.99'
$ws.Range("E39").Value = '  -4.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01609'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.83%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.889'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.85%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8604'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.004'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.77'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("D45").Value = '1.840.25'
$ws.Range("E45").Value = '  -5.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000113'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.043'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4312'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.61%  '
$ws.Range("E51").Value = '  -4.40%  '
